$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on D-column cells whose new values would otherwise be
# auto-parsed by Excel as numbers, so they stay as literal text like the source data.
$textForceRows = @(4,5,6,7,8,9,10,11,12,14,15,16,18,20,22,23,24,25,26,27,28,29,30,31,32,33,34,36,37,38,39,40,41,42,43,44,45,46,47,48,51)
foreach ($r in $textForceRows) {
    $ws.Range("D" + $r).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "30.469.74"
$ws.Range("E2").Value = "  -1.07%  "

# Row 3
$ws.Range("D3").Value = "1.923.08"
$ws.Range("E3").Value = "  +1.79%  "

# Row 4
$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").Value = "  -0.13%  "

# Row 5
$ws.Range("D5").Value = "243.01"
$ws.Range("E5").Value = "  +1.27%  "

# Row 6
$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").Value = "  -0.16%  "

# Row 7
$ws.Range("D7").Value = "0.4697"
$ws.Range("E7").Value = "  -1.82%  "

# Row 8
$ws.Range("D8").Value = "0.2880"
$ws.Range("E8").Value = "  -2.70%  "

# Row 9
$ws.Range("D9").Value = "0.06799"
$ws.Range("E9").Value = "  +2.51%  "

# Row 10
$ws.Range("D10").Value = "109.22"
$ws.Range("E10").Value = "  +8.48%  "

# Row 11
$ws.Range("D11").Value = "18.38"
$ws.Range("E11").Value = "  -1.70%  "

# Row 12
$ws.Range("D12").Value = "0.07730"
$ws.Range("E12").Value = "  +2.37%  "

# Row 13
$ws.Range("D13").Value = "1.890.49"
$ws.Range("E13").Value = "  +0.36%  "

# Row 14
$ws.Range("D14").Value = "5.320"
$ws.Range("E14").Value = "  +3.29%  "

# Row 15
$ws.Range("D15").Value = "0.6598"
$ws.Range("E15").Value = "  -0.27%  "

# Row 16
$ws.Range("D16").Value = "295.99"
$ws.Range("E16").Value = "  -2.21%  "

# Row 17
$ws.Range("D17").Value = "30.464.30"
$ws.Range("E17").Value = "  -1.05%  "

# Row 18
$ws.Range("D18").Value = "0.000007622"
$ws.Range("E18").Value = "  +0.40%  "

# Row 19
$ws.Range("E19").Value = "  -1.23%  "

# Row 20
$ws.Range("D20").Value = "0.9997"
$ws.Range("E20").Value = "  -0.07%  "

# Row 21
$ws.Range("D21").Value = "2.139.87"
$ws.Range("E21").Value = "  +0.31%  "

# Row 22
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.15%  "

# Row 23
$ws.Range("D23").Value = "5.251"
$ws.Range("E23").Value = "  +1.83%  "

# Row 24
$ws.Range("D24").Value = "6.206"
$ws.Range("E24").Value = "  -0.21%  "

# Row 25
$ws.Range("D25").Value = "9.390"
$ws.Range("E25").Value = "  +0.69%  "

# Row 26
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "168.79"
$ws.Range("E26").Value = "  +0.55%  "

# Row 27
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "21.61"
$ws.Range("E27").Value = "  +5.91%  "

# Row 28
$ws.Range("D28").Value = "2.130"
$ws.Range("E28").Value = "  +8.96%  "

# Row 29
$ws.Range("D29").Value = "0.1069"
$ws.Range("E29").Value = "  -5.23%  "

# Row 30
$ws.Range("D30").Value = "1.365"
$ws.Range("E30").Value = "  +0.91%  "

# Row 31
$ws.Range("D31").Value = "4.189"
$ws.Range("E31").Value = "  +0.35%  "

# Row 32
$ws.Range("D32").Value = "3.996"
$ws.Range("E32").Value = "  -0.05%  "

# Row 33
$ws.Range("D33").Value = "0.05043"
$ws.Range("E33").Value = "  -0.86%  "

# Row 34
$ws.Range("D34").Value = "0.7399"
$ws.Range("E34").Value = "  +0.82%  "

# Row 35
$ws.Range("E35").Value = "  -0.72%  "

# Row 36
$ws.Range("D36").Value = "0.02100"
$ws.Range("E36").Value = "  +6.30%  "

# Row 37
$ws.Range("D37").Value = "2.738"
$ws.Range("E37").Value = "  +0.75%  "

# Row 38
$ws.Range("D38").Value = "2.684"
$ws.Range("E38").Value = "  -0.83%  "

# Row 39
$ws.Range("D39").Value = "2.067"
$ws.Range("E39").Value = "  +0.29%  "

# Row 40
$ws.Range("D40").Value = "110.25"
$ws.Range("E40").Value = "  +2.02%  "

# Row 41
$ws.Range("D41").Value = "0.8721"
$ws.Range("E41").Value = "  -2.78%  "

# Row 42
$ws.Range("D42").Value = "5.871"
$ws.Range("E42").Value = "  +3.83%  "

# Row 43
$ws.Range("D43").Value = "0.4254"
$ws.Range("E43").Value = "  +1.25%  "

# Row 44
$ws.Range("D44").Value = "0.9996"
$ws.Range("E44").Value = "  -0.18%  "

# Row 45
$ws.Range("D45").Value = "67.51"
$ws.Range("E45").Value = "  +0.82%  "

# Row 46
$ws.Range("D46").Value = "51.10"
$ws.Range("E46").Value = "  +19.17%  "

# Row 47
$ws.Range("D47").Value = "7.197"
$ws.Range("E47").Value = "  -2.38%  "

# Row 48
$ws.Range("D48").Value = "9.298"
$ws.Range("E48").Value = "  +1.84%  "

# Row 49
$ws.Range("E49").Value = "  -0.94%  "

# Row 50
$ws.Range("E50").Value = "  +0.46%  "

# Row 51
$ws.Range("D51").Value = "0.2488"
$ws.Range("E51").Value = "  +11.33%  "
